$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("iciba-collins")

# Duplicate column C (Define) formatting/width into a new column, inserted
# before the old column D (QC) -- shifts QC/QT right to E/F.
$ws.Columns.Item(3).Copy()
$ws.Columns.Item(4).Insert()

# New column D header
$ws.Range("D1").Value = "Note"

$ws.AutoFilter.Range.Select()
